# Making searchable recipes based on ingredients
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shorten source names
$ws.Range("B2").Value = "Claire B"
$ws.Range("B3").Value = "Claire C"
$ws.Range("B6").Value = "Mitchell B"
$ws.Range("B7").Value = "Mum B"
$ws.Range("B8").Value = "Mum B"
$ws.Range("B9").Value = "Shalane F"
$ws.Range("B10").Value = "Shalane F"
$ws.Range("B11").Value = "Shalane F"
$ws.Range("B12").Value = "Shalane F"
$ws.Range("B13").Value = "Shalane F"
$ws.Range("B14").Value = "Shalane F"
$ws.Range("B15").Value = "Shalane F"
$ws.Range("B16").Value = "Shalane F"

# Fill in TBD directions for recipes missing them (searchable placeholder)
$ws.Range("C8").Value = "TBD"
$ws.Range("C16").Value = "TBD"

# Row height tweak on the Chicken Noodle Soup row (wrapped text row)
$ws.Rows.Item(6).RowHeight = 23.85

# Update selection to reflect last edited cell
$ws.Range("A16").Select()
